# Fill in the two previously-missing ("na") values for row 12 (person with
# missing Age / HS Graduate) on the "Data" sheet, and update the selected
# cell to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 12: Age was "na" -> now a real numeric age of 50
$ws.Range("D12").Value = 50

# Row 12: HS Graduate was "na" -> now "Yes"
$ws.Range("E12").Value = "Yes"

# Update the active selection on the Data sheet to K11 (matches saved view)
$ws.Activate()
$ws.Range("K11").Select()
